# Weekly fruit/veg price update: insert a new observation row for
# "Acelga" (Agrícola del Norte S.A. de Arica) at row 46, pushing the
# existing rows 46:66 down to 47:67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46 (existing rows shift down by one).
$ws.Rows(46).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A46").Value = 1
$ws.Range("B46").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C46").Value = "Arica y Parinacota"
$ws.Range("D46").Value = 44755
$ws.Range("E46").Value = 15
$ws.Range("F46").Value = 100112009
$ws.Range("G46").Value = "Acelga"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 250
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = 2250
$ws.Range("N46").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value = 750
$ws.Range("Q46").Value = 3
$ws.Range("R46").Value = "Hortaliza"
